$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# fm21 row (row 26) testing is now complete
$ws.Range("I26").Value = "complete"
$ws.Range("J26").Value = "complete"

# Add new row 27 for fm22, based on the formatting of row 26
$ws.Range("B26:J26").Copy($ws.Range("B27:J27"))

$ws.Range("B27").Value = "fm22"
$ws.Range("C27").Value = "Special condition sublimits, blanket policy deductible and multiple policy layers, with IL back-allocation"
$ws.Range("D27").Value = "0,2"
$ws.Range("E27").Value = "2,11,12,14"
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = "3,5,6,12"
$ws.Range("I27").Value = "complete"
$ws.Range("J27").Value = "complete"

# D27/E27 should have the default (no) style, matching the source diff
$ws.Range("D27:E27").Style = "Normal"

# Update the selected cell shown in the sheet view
[void]$ws.Range("A26").Select()
